# Update the "Training Dashboard" sheet: for rows 3-11 decrement the
# "PERIOD TO EXPIRE" (column H) counts by 1 day and bump the
# "LAST UPDATE" (column I) date text from 03-Nov-2025 to 04-Nov-2025,
# reflecting one more day of progress as of 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 11; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE
    # Assign via a text-returning formula first so Excel does not
    # auto-detect the "dd-mmm-yyyy" literal as a date (which would
    # otherwise flip the cell to a date serial + date number format and
    # mint a new style). Then convert the formula result down to a
    # plain literal value, keeping the original style/format intact.
    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy()
    $iCell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
